$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("BOM")

# Row 31: Fuse Holder SMD
$ws.Range("A31").Value = "Fuse Holder SMD"
$ws.Range("A31").Style = $ws.Range("A30").Style
$ws.Range("C31").Value = "F1-F5"
$ws.Range("G31").Value = 4.0599999999999996
$ws.Range("G31").Style = $ws.Range("G30").Style
$ws.Range("I31").Value = 5
$ws.Range("I31").Style = $ws.Range("I30").Style
$ws.Range("K31").Formula = "=I31*G31"
$ws.Range("M31").Value = "https://www.digikey.ca/en/products/detail/littelfuse-inc/0154003-DR/183356"

# Row 32: Fuse 2-SMD
$ws.Range("A32").Value = "Fuse 2-SMD"
$ws.Range("A32").Style = $ws.Range("A30").Style
$ws.Range("C32").Value = "F1-F5"

$ws.Range("A4").Select()
